# Daily COVID data upload: correct the running index in column A for the
# most recent date block on each sheet (it had incorrectly restarted at 0),
# then append the newly-published 2020-12-14 rows to both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: cases_by_race
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("cases_by_race")

# Fix the running counter in column A for rows 44-52 (2020-12-13 block):
# it had reset to 0-8, it should continue on from the prior block (42-50).
$raceFix = @{ 44 = 42; 45 = 43; 46 = 44; 47 = 45; 48 = 46; 49 = 47; 50 = 48; 51 = 49; 52 = 50 }
foreach ($r in $raceFix.Keys) {
    $ws1.Cells.Item($r, 1).Value = $raceFix[$r]
}

# New rows for the 2020-12-14 refresh (DATA_AS_OF_DT 2020-12-13), appended
# after row 52. Each row's format is copied from the matching row in the
# prior (2020-12-13) block so styles/borders on column A are preserved.
$raceNewRows = @(
    @{ Row = 53; Src = 44; A = 0; B = "";                                           E = 4 },
    @{ Row = 54; Src = 45; A = 1; B = "American Indian or Alaska Native";           E = 49 },
    @{ Row = 55; Src = 46; A = 2; B = "Asian";                                      E = 231 },
    @{ Row = 56; Src = 47; A = 3; B = "Black or African American";                  E = 1354 },
    @{ Row = 57; Src = 48; A = 4; B = "Native Hawaiian or Other Pacific Islander";   E = 10 },
    @{ Row = 58; Src = 49; A = 5; B = "Not disclosed";                              E = 1514 },
    @{ Row = 59; Src = 50; A = 6; B = "Other Race";                                 E = 363 },
    @{ Row = 60; Src = 51; A = 7; B = "Two or more";                                E = 131 },
    @{ Row = 61; Src = 52; A = 8; B = "White";                                      E = 12693 }
)

foreach ($item in $raceNewRows) {
    $srcRange = $ws1.Range("A" + $item.Src + ":E" + $item.Src)
    $dstRange = $ws1.Range("A" + $item.Row + ":E" + $item.Row)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats

    # Keep the date columns as literal text (not auto-converted to dates).
    $ws1.Range("C" + $item.Row + ":D" + $item.Row).NumberFormat = "@"

    $ws1.Cells.Item($item.Row, 1).Value = $item.A
    $ws1.Cells.Item($item.Row, 2).Value = $item.B
    $ws1.Cells.Item($item.Row, 3).Value = "2020-12-14"
    $ws1.Cells.Item($item.Row, 4).Value = "2020-12-13"
    $ws1.Cells.Item($item.Row, 5).Value = $item.E
}

# ---------------------------------------------------------------------
# Sheet 2: cases_by_ethnicity
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("cases_by_ethnicity")

# Fix the running counter in column A for rows 20-22 (2020-12-13 block).
$ethFix = @{ 20 = 18; 21 = 19; 22 = 20 }
foreach ($r in $ethFix.Keys) {
    $ws2.Cells.Item($r, 1).Value = $ethFix[$r]
}

# New rows for the 2020-12-14 refresh.
$ethNewRows = @(
    @{ Row = 23; Src = 20; A = 0; B = "Hispanic or Latino";     E = 361 },
    @{ Row = 24; Src = 21; A = 1; B = "Not Hispanic or Latino"; E = 12722 },
    @{ Row = 25; Src = 22; A = 2; B = "unknown";                E = 3266 }
)

foreach ($item in $ethNewRows) {
    $srcRange = $ws2.Range("A" + $item.Src + ":E" + $item.Src)
    $dstRange = $ws2.Range("A" + $item.Row + ":E" + $item.Row)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats

    $ws2.Range("C" + $item.Row + ":D" + $item.Row).NumberFormat = "@"

    $ws2.Cells.Item($item.Row, 1).Value = $item.A
    $ws2.Cells.Item($item.Row, 2).Value = $item.B
    $ws2.Cells.Item($item.Row, 3).Value = "2020-12-14"
    $ws2.Cells.Item($item.Row, 4).Value = "2020-12-13"
    $ws2.Cells.Item($item.Row, 5).Value = $item.E
}
